# Insert a new record at row 101 (pushing the existing rows 101-146 down to
# 102-147) for "Terminal La Palmera de La Serena" / Ajo / Chino / Primera,
# dated 2021-09-27 (serial 44466).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 101..146 down to 102..147, duplicating formatting (incl. the
# date style on column D) from the row being pushed down, just like using
# Excel's "Insert" on a selected row.
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new observation.
$ws.Cells.Item(101, 1).Value  = 8
$ws.Cells.Item(101, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(101, 3).Value  = "Coquimbo"
$ws.Cells.Item(101, 4).Value  = 44466
$ws.Cells.Item(101, 5).Value  = 4
$ws.Cells.Item(101, 6).Value  = 100112003
$ws.Cells.Item(101, 7).Value  = "Ajo"
$ws.Cells.Item(101, 8).Value  = "Chino"
$ws.Cells.Item(101, 9).Value  = "Primera"
$ws.Cells.Item(101, 10).Value = 640
$ws.Cells.Item(101, 11).Value = 16000
$ws.Cells.Item(101, 12).Value = 17000
$ws.Cells.Item(101, 13).Value = 16500
$ws.Cells.Item(101, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(101, 15).Value = "China"
$ws.Cells.Item(101, 16).Value = 1650
$ws.Cells.Item(101, 17).Value = 10
$ws.Cells.Item(101, 18).Value = "Hortaliza"
